# Add the new Q&A row (ID=2, Question="Is the grant available to a startup?")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Is the grant available to a startup?"

# Narrow column A (was sized for the long header text)
$ws.Columns.Item(1).ColumnWidth = 9.26

# Leave the cursor where the author's session ended up
$ws.Range("C7").Select()
